$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-11-16 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-11-17 Monday", 2)

# Update the division problems in the table. Cells are addressed directly
# by (row, column) and their Range.Text is replaced in place, so that
# duplicate "before" values (e.g. the two "25÷8=" cells) each end up with
# their own distinct "after" value rather than a document-wide Find
# touching the wrong occurrence.
$t = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; New = "62÷8=" },
    @{ Row = 1;  Col = 2; New = "64÷3=" },
    @{ Row = 1;  Col = 3; New = "22÷3=" },
    @{ Row = 1;  Col = 4; New = "38÷5=" },
    @{ Row = 1;  Col = 5; New = "33÷9=" },

    @{ Row = 5;  Col = 1; New = "24÷9=" },
    @{ Row = 5;  Col = 2; New = "94÷6=" },
    @{ Row = 5;  Col = 3; New = "48÷4=" },
    @{ Row = 5;  Col = 4; New = "56÷9=" },
    @{ Row = 5;  Col = 5; New = "42÷5=" },

    @{ Row = 9;  Col = 1; New = "36÷2=" },
    @{ Row = 9;  Col = 2; New = "90÷9=" },
    @{ Row = 9;  Col = 3; New = "18÷2=" },
    @{ Row = 9;  Col = 4; New = "47÷9=" },
    @{ Row = 9;  Col = 5; New = "80÷9=" },

    @{ Row = 13; Col = 1; New = "83÷5=" },
    @{ Row = 13; Col = 2; New = "40÷8=" },
    @{ Row = 13; Col = 3; New = "27÷6=" },
    @{ Row = 13; Col = 4; New = "12÷9=" },
    @{ Row = 13; Col = 5; New = "35÷2=" },

    @{ Row = 17; Col = 1; New = "69÷4=" },
    @{ Row = 17; Col = 2; New = "33÷7=" },
    @{ Row = 17; Col = 3; New = "75÷6=" },
    @{ Row = 17; Col = 4; New = "61÷4=" },
    @{ Row = 17; Col = 5; New = "82÷4=" }
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $cell.Range.Text = $item.New
}
